$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update validation / test metric values ---
$ws.Range("F4").Value2 = 70.09
$ws.Range("G4").Value2 = 79.7
$ws.Range("G5").Value2 = 58.83
$ws.Range("G6").Value2 = 74.56
$ws.Range("F7").Value2 = 68.69
$ws.Range("G7").Value2 = 70.28
$ws.Range("G10").Value2 = 75.84
$ws.Range("F12").Value2 = 74.07
$ws.Range("F13").Value2 = 72.34
$ws.Range("E15").Value2 = 81.17
$ws.Range("F15").Value2 = 77.03
$ws.Range("G15").Value2 = 74.72
$ws.Range("E16").Value2 = 79.44
$ws.Range("F16").Value2 = 75.73
$ws.Range("G16").Value2 = 81.7
$ws.Range("E17").Value2 = 79.26
$ws.Range("F17").Value2 = 75.97
$ws.Range("G17").Value2 = 79.86
$ws.Range("E18").Value2 = 79.03
$ws.Range("F18").Value2 = 76.47
$ws.Range("G18").Value2 = 82.89
$ws.Range("E19").Value2 = 75.1
$ws.Range("F19").Value2 = 72.9
$ws.Range("G19").Value2 = 87.7
$ws.Range("E20").Value2 = 81.03
$ws.Range("F20").Value2 = 76.92
$ws.Range("G20").Value2 = 81.54

# --- Fix borders for rows 14-20 (threshold table) ---
# E14
$ws.Range("E14").Borders.Item(7).LineStyle = 1
$ws.Range("E14").Borders.Item(7).Weight = 2
$ws.Range("E14").Borders.Item(8).LineStyle = -4142
$ws.Range("E14").Borders.Item(9).LineStyle = 1
$ws.Range("E14").Borders.Item(9).Weight = -4138
$ws.Range("E14").Borders.Item(10).LineStyle = 1
$ws.Range("E14").Borders.Item(10).Weight = 2
# F14
$ws.Range("F14").Borders.Item(7).LineStyle = 1
$ws.Range("F14").Borders.Item(7).Weight = 2
$ws.Range("F14").Borders.Item(8).LineStyle = -4142
$ws.Range("F14").Borders.Item(9).LineStyle = 1
$ws.Range("F14").Borders.Item(9).Weight = -4138
$ws.Range("F14").Borders.Item(10).LineStyle = 1
$ws.Range("F14").Borders.Item(10).Weight = 2
# G14
$ws.Range("G14").Borders.Item(7).LineStyle = -4142
$ws.Range("G14").Borders.Item(8).LineStyle = -4142
$ws.Range("G14").Borders.Item(9).LineStyle = 1
$ws.Range("G14").Borders.Item(9).Weight = -4138
$ws.Range("G14").Borders.Item(10).LineStyle = 1
$ws.Range("G14").Borders.Item(10).Weight = -4138
# C15
$ws.Range("C15").Borders.Item(7).LineStyle = 1
$ws.Range("C15").Borders.Item(7).Weight = 2
$ws.Range("C15").Borders.Item(8).LineStyle = 1
$ws.Range("C15").Borders.Item(8).Weight = -4138
$ws.Range("C15").Borders.Item(9).LineStyle = -4142
$ws.Range("C15").Borders.Item(10).LineStyle = -4142
# D15
$ws.Range("D15").Borders.Item(7).LineStyle = 1
$ws.Range("D15").Borders.Item(7).Weight = 2
$ws.Range("D15").Borders.Item(8).LineStyle = -4142
$ws.Range("D15").Borders.Item(9).LineStyle = -4142
$ws.Range("D15").Borders.Item(10).LineStyle = 1
$ws.Range("D15").Borders.Item(10).Weight = 2
# E15
$ws.Range("E15").Borders.Item(7).LineStyle = -4142
$ws.Range("E15").Borders.Item(8).LineStyle = -4142
$ws.Range("E15").Borders.Item(9).LineStyle = -4142
$ws.Range("E15").Borders.Item(10).LineStyle = -4142
# F15
$ws.Range("F15").Borders.Item(7).LineStyle = 1
$ws.Range("F15").Borders.Item(7).Weight = 2
$ws.Range("F15").Borders.Item(8).LineStyle = -4142
$ws.Range("F15").Borders.Item(9).LineStyle = -4142
$ws.Range("F15").Borders.Item(10).LineStyle = 1
$ws.Range("F15").Borders.Item(10).Weight = 2
# G15
$ws.Range("G15").Borders.Item(7).LineStyle = 1
$ws.Range("G15").Borders.Item(7).Weight = 2
$ws.Range("G15").Borders.Item(8).LineStyle = -4142
$ws.Range("G15").Borders.Item(9).LineStyle = -4142
$ws.Range("G15").Borders.Item(10).LineStyle = 1
$ws.Range("G15").Borders.Item(10).Weight = -4138
# C16
$ws.Range("C16").Borders.Item(7).LineStyle = 1
$ws.Range("C16").Borders.Item(7).Weight = 2
$ws.Range("C16").Borders.Item(8).LineStyle = -4142
$ws.Range("C16").Borders.Item(9).LineStyle = -4142
$ws.Range("C16").Borders.Item(10).LineStyle = -4142
# D16
$ws.Range("D16").Borders.Item(7).LineStyle = 1
$ws.Range("D16").Borders.Item(7).Weight = 2
$ws.Range("D16").Borders.Item(8).LineStyle = -4142
$ws.Range("D16").Borders.Item(9).LineStyle = -4142
$ws.Range("D16").Borders.Item(10).LineStyle = 1
$ws.Range("D16").Borders.Item(10).Weight = 2
# E16
$ws.Range("E16").Borders.Item(7).LineStyle = -4142
$ws.Range("E16").Borders.Item(8).LineStyle = -4142
$ws.Range("E16").Borders.Item(9).LineStyle = -4142
$ws.Range("E16").Borders.Item(10).LineStyle = -4142
# F16
$ws.Range("F16").Borders.Item(7).LineStyle = 1
$ws.Range("F16").Borders.Item(7).Weight = 2
$ws.Range("F16").Borders.Item(8).LineStyle = -4142
$ws.Range("F16").Borders.Item(9).LineStyle = -4142
$ws.Range("F16").Borders.Item(10).LineStyle = 1
$ws.Range("F16").Borders.Item(10).Weight = 2
# G16
$ws.Range("G16").Borders.Item(7).LineStyle = 1
$ws.Range("G16").Borders.Item(7).Weight = 2
$ws.Range("G16").Borders.Item(8).LineStyle = -4142
$ws.Range("G16").Borders.Item(9).LineStyle = -4142
$ws.Range("G16").Borders.Item(10).LineStyle = 1
$ws.Range("G16").Borders.Item(10).Weight = -4138
# C17
$ws.Range("C17").Borders.Item(7).LineStyle = 1
$ws.Range("C17").Borders.Item(7).Weight = 2
$ws.Range("C17").Borders.Item(8).LineStyle = -4142
$ws.Range("C17").Borders.Item(9).LineStyle = 1
$ws.Range("C17").Borders.Item(9).Weight = 2
$ws.Range("C17").Borders.Item(10).LineStyle = -4142
# D17
$ws.Range("D17").Borders.Item(7).LineStyle = 1
$ws.Range("D17").Borders.Item(7).Weight = 2
$ws.Range("D17").Borders.Item(8).LineStyle = -4142
$ws.Range("D17").Borders.Item(9).LineStyle = 1
$ws.Range("D17").Borders.Item(9).Weight = 2
$ws.Range("D17").Borders.Item(10).LineStyle = 1
$ws.Range("D17").Borders.Item(10).Weight = 2
# E17
$ws.Range("E17").Borders.Item(7).LineStyle = -4142
$ws.Range("E17").Borders.Item(8).LineStyle = -4142
$ws.Range("E17").Borders.Item(9).LineStyle = -4142
$ws.Range("E17").Borders.Item(10).LineStyle = -4142
# F17
$ws.Range("F17").Borders.Item(7).LineStyle = 1
$ws.Range("F17").Borders.Item(7).Weight = 2
$ws.Range("F17").Borders.Item(8).LineStyle = -4142
$ws.Range("F17").Borders.Item(9).LineStyle = -4142
$ws.Range("F17").Borders.Item(10).LineStyle = 1
$ws.Range("F17").Borders.Item(10).Weight = 2
# G17
$ws.Range("G17").Borders.Item(7).LineStyle = 1
$ws.Range("G17").Borders.Item(7).Weight = 2
$ws.Range("G17").Borders.Item(8).LineStyle = -4142
$ws.Range("G17").Borders.Item(9).LineStyle = -4142
$ws.Range("G17").Borders.Item(10).LineStyle = 1
$ws.Range("G17").Borders.Item(10).Weight = -4138
# C18
$ws.Range("C18").Borders.Item(7).LineStyle = 1
$ws.Range("C18").Borders.Item(7).Weight = 2
$ws.Range("C18").Borders.Item(8).LineStyle = 1
$ws.Range("C18").Borders.Item(8).Weight = 2
$ws.Range("C18").Borders.Item(9).LineStyle = -4142
$ws.Range("C18").Borders.Item(10).LineStyle = -4142
# D18
$ws.Range("D18").Borders.Item(7).LineStyle = 1
$ws.Range("D18").Borders.Item(7).Weight = 2
$ws.Range("D18").Borders.Item(8).LineStyle = -4142
$ws.Range("D18").Borders.Item(9).LineStyle = -4142
$ws.Range("D18").Borders.Item(10).LineStyle = 1
$ws.Range("D18").Borders.Item(10).Weight = 2
# E18
$ws.Range("E18").Borders.Item(7).LineStyle = -4142
$ws.Range("E18").Borders.Item(8).LineStyle = 1
$ws.Range("E18").Borders.Item(8).Weight = 2
$ws.Range("E18").Borders.Item(9).LineStyle = -4142
$ws.Range("E18").Borders.Item(10).LineStyle = -4142
# F18
$ws.Range("F18").Borders.Item(7).LineStyle = 1
$ws.Range("F18").Borders.Item(7).Weight = 2
$ws.Range("F18").Borders.Item(8).LineStyle = 1
$ws.Range("F18").Borders.Item(8).Weight = 2
$ws.Range("F18").Borders.Item(9).LineStyle = -4142
$ws.Range("F18").Borders.Item(10).LineStyle = -4142
# G18
$ws.Range("G18").Borders.Item(7).LineStyle = 1
$ws.Range("G18").Borders.Item(7).Weight = 2
$ws.Range("G18").Borders.Item(8).LineStyle = 1
$ws.Range("G18").Borders.Item(8).Weight = 2
$ws.Range("G18").Borders.Item(9).LineStyle = -4142
$ws.Range("G18").Borders.Item(10).LineStyle = 1
$ws.Range("G18").Borders.Item(10).Weight = -4138
# C19
$ws.Range("C19").Borders.Item(7).LineStyle = 1
$ws.Range("C19").Borders.Item(7).Weight = 2
$ws.Range("C19").Borders.Item(8).LineStyle = -4142
$ws.Range("C19").Borders.Item(9).LineStyle = -4142
$ws.Range("C19").Borders.Item(10).LineStyle = -4142
# D19
$ws.Range("D19").Borders.Item(7).LineStyle = 1
$ws.Range("D19").Borders.Item(7).Weight = 2
$ws.Range("D19").Borders.Item(8).LineStyle = -4142
$ws.Range("D19").Borders.Item(9).LineStyle = -4142
$ws.Range("D19").Borders.Item(10).LineStyle = 1
$ws.Range("D19").Borders.Item(10).Weight = 2
# E19
$ws.Range("E19").Borders.Item(7).LineStyle = -4142
$ws.Range("E19").Borders.Item(8).LineStyle = -4142
$ws.Range("E19").Borders.Item(9).LineStyle = -4142
$ws.Range("E19").Borders.Item(10).LineStyle = -4142
# F19
$ws.Range("F19").Borders.Item(7).LineStyle = 1
$ws.Range("F19").Borders.Item(7).Weight = 2
$ws.Range("F19").Borders.Item(8).LineStyle = -4142
$ws.Range("F19").Borders.Item(9).LineStyle = -4142
$ws.Range("F19").Borders.Item(10).LineStyle = -4142
# G19
$ws.Range("G19").Borders.Item(7).LineStyle = 1
$ws.Range("G19").Borders.Item(7).Weight = 2
$ws.Range("G19").Borders.Item(8).LineStyle = -4142
$ws.Range("G19").Borders.Item(9).LineStyle = -4142
$ws.Range("G19").Borders.Item(10).LineStyle = 1
$ws.Range("G19").Borders.Item(10).Weight = -4138
# D20
$ws.Range("D20").Borders.Item(7).LineStyle = 1
$ws.Range("D20").Borders.Item(7).Weight = 2
$ws.Range("D20").Borders.Item(8).LineStyle = -4142
$ws.Range("D20").Borders.Item(9).LineStyle = 1
$ws.Range("D20").Borders.Item(9).Weight = -4138
$ws.Range("D20").Borders.Item(10).LineStyle = -4142
# E20
$ws.Range("E20").Borders.Item(7).LineStyle = 1
$ws.Range("E20").Borders.Item(7).Weight = 2
$ws.Range("E20").Borders.Item(8).LineStyle = -4142
$ws.Range("E20").Borders.Item(9).LineStyle = 1
$ws.Range("E20").Borders.Item(9).Weight = -4138
$ws.Range("E20").Borders.Item(10).LineStyle = -4142
# F20
$ws.Range("F20").Borders.Item(7).LineStyle = 1
$ws.Range("F20").Borders.Item(7).Weight = 2
$ws.Range("F20").Borders.Item(8).LineStyle = -4142
$ws.Range("F20").Borders.Item(9).LineStyle = 1
$ws.Range("F20").Borders.Item(9).Weight = -4138
$ws.Range("F20").Borders.Item(10).LineStyle = -4142
# G20
$ws.Range("G20").Borders.Item(7).LineStyle = 1
$ws.Range("G20").Borders.Item(7).Weight = 2
$ws.Range("G20").Borders.Item(8).LineStyle = -4142
$ws.Range("G20").Borders.Item(9).LineStyle = 1
$ws.Range("G20").Borders.Item(9).Weight = -4138
$ws.Range("G20").Borders.Item(10).LineStyle = 1
$ws.Range("G20").Borders.Item(10).Weight = -4138

# --- Update selection ---
$ws.Range("F12").Select()
